# Updated cryptos list (price + 1h volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "65.665.17"
$cell.ClearFormats()
$ws.Range("E2").Value = "  -0.05%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.667.55"
$cell.ClearFormats()
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.03%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "599.30"
$cell.ClearFormats()
$ws.Range("E5").Value = "  -1.19%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "155.45"
$cell.ClearFormats()
$ws.Range("E6").Value = "  -1.92%  "

$ws.Range("E7").Value = "  +0.08%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.602"
$cell.ClearFormats()
$ws.Range("E8").Value = "  +1.91%  "

$ws.Range("E9").Value = "  -1.55%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "5.90"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +0.97%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.396"
$cell.ClearFormats()
$ws.Range("E11").Value = "  -2.36%  "

$ws.Range("E12").Value = "  -0.21%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "29.27"
$cell.ClearFormats()
$ws.Range("E13").Value = "  -1.92%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "0.0000194"
$cell.ClearFormats()
$ws.Range("E14").Value = "  +0.13%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "3.153.64"
$cell.ClearFormats()
$ws.Range("E15").Value = "  +0.15%  "

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "65.485.05"
$cell.ClearFormats()
$ws.Range("E16").Value = "  +0.06%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.679.09"
$cell.ClearFormats()
$ws.Range("E17").Value = "  +0.37%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "12.50"
$cell.ClearFormats()
$ws.Range("E18").Value = "  -2.22%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "4.80"
$cell.ClearFormats()
$ws.Range("E19").Value = "  -2.22%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.48"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +1.32%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "350.01"
$cell.ClearFormats()
$ws.Range("E21").Value = "  -2.95%  "

$ws.Range("E22").Value = "  -0.15%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "70.10"
$cell.ClearFormats()
$ws.Range("E23").Value = "  +1.47%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.73"
$cell.ClearFormats()
$ws.Range("E24").Value = "  +1.55%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "0.0000108"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +2.26%  "

$ws.Range("E26").Value = "  -3.82%  "

$ws.Range("E27").Value = "  -1.94%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "0.168"
$cell.ClearFormats()
$ws.Range("E28").Value = "  +1.37%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "8.08"
$cell.ClearFormats()
$ws.Range("E29").Value = "  -1.91%  "

$ws.Range("E30").Value = "  -0.44%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "535.69"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -1.30%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.ClearFormats()
$ws.Range("E32").Value = "  -2.52%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.75"
$cell.ClearFormats()
$ws.Range("E33").Value = "  -5.50%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.52"
$cell.ClearFormats()
$ws.Range("E34").Value = "  +2.29%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.39"
$cell.ClearFormats()
$ws.Range("E35").Value = "  -4.59%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.423"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -2.74%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "20.33"
$cell.ClearFormats()
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("E39").Value = "  -0.01%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "1.94"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -4.07%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "42.43"
$cell.ClearFormats()
$ws.Range("E42").Value = "  -0.04%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "165.88"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.67%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "4.07"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -2.88%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.0608"
$cell.ClearFormats()
$ws.Range("E45").Value = "  -1.06%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "22.92"
$cell.ClearFormats()
$ws.Range("E46").Value = "  -1.37%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.23"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -6.19%  "

$ws.Range("E48").Value = "  -2.52%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.0259"
$cell.ClearFormats()
$ws.Range("E49").Value = "  -2.53%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.0995"
$cell.ClearFormats()
$ws.Range("E50").Value = "  +0.32%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "19.90"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -0.06%  "

